$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '87.207.19'
$ws.Range("E2").Value = '  -3.58%  '
$ws.Range("D3").Value = '3.053.58'
$ws.Range("E3").Value = '  -4.16%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '621.52'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.361'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -9.39%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.773'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +11.91%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("D10").Value = '3.051.98'
$ws.Range("E10").Value = '  -4.12%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.581'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("E12").Value = '  -0.25%  '
$ws.Range("E13").Value = '  -9.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.23'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.49%  '
$ws.Range("D15").Value = '86.943.73'
$ws.Range("E15").Value = '  -3.53%  '
$ws.Range("D16").Value = '3.607.96'
$ws.Range("E16").Value = '  -4.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '31.18'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.77%  '
$ws.Range("D18").Value = '3.074.43'
$ws.Range("E18").Value = '  -3.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.35'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000206'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.97'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '414.32'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.77'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.43'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '81.88'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +8.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.23'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.98%  '
$ws.Range("D28").Value = '3.217.46'
$ws.Range("E28").Value = '  -4.13%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.148'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -13.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.98'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '494.79'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.59'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -13.91%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.140'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +9.90%  '
$ws.Range("E36").Value = '  -6.17%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.78'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.04%  '
$ws.Range("E38").Value = '  -2.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '21.78'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.44%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.14'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.06%  '
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("E43").Value = '  -4.60%  '
$ws.Range("E44").Value = '  -7.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '146.59'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.131'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.22%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.38'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0638'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '157.81'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.65%  '
$ws.Range("E50").Value = '  -1.13%  '
$ws.Range("E51").Value = '  -6.25%  '
